$d = $word.ActiveDocument

function Get-ParagraphIndexContaining($doc, $needle) {
    $count = $doc.Paragraphs.Count
    for ($i = 1; $i -le $count; $i++) {
        $p = $doc.Paragraphs.Item($i)
        if ($p.Range.Text -like ("*" + $needle + "*")) {
            return $i
        }
    }
    return -1
}

# ---------------------------------------------------------------------------
# Change 1: add a new participant line for "Xing Meng (u6483085)" right after
# the existing "Yafei Liu (u6605935)" participant paragraph.
# ---------------------------------------------------------------------------
$yafeiIdx = Get-ParagraphIndexContaining $d "u6605935"
$yafeiPara = $d.Paragraphs.Item($yafeiIdx)
$yafeiRange = $d.Range($yafeiPara.Range.Start, $yafeiPara.Range.End)
$yafeiRange.InsertParagraphAfter()

$newPara = $d.Paragraphs.Item($yafeiIdx + 1)
$newRange = $d.Range($newPara.Range.Start, $newPara.Range.End)

$newParaXml = '<?xml version="1.0" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r><w:tab/><w:t xml:space="preserve">         Xing Meng (u6483085)</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$null = $newRange.InsertXML($newParaXml)

# ---------------------------------------------------------------------------
# Change 2: mark "Yafei" as a flagged-spelling word (spellStart/spellEnd) in
# the "Yafei Liu, Jiale Wang - Build Website" task-assignment line, splitting
# the leading run exactly as Word's proofing pass would, while leaving the
# rest of that paragraph (Jiale's existing proofErr wrapper, etc.) intact.
# The whole paragraph is rewritten in one InsertXML call (on its own full
# Range) so the untouched tail content is carried over byte-for-byte and the
# paragraph's own identity (paraId/rsid/pPr) is preserved.
# ---------------------------------------------------------------------------
$taskIdx = Get-ParagraphIndexContaining $d "Build Website"
$taskPara = $d.Paragraphs.Item($taskIdx)
$taskRange = $d.Range($taskPara.Range.Start, $taskPara.Range.End)

$splitXml = '<?xml version="1.0" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml"><w:body><w:p w14:paraId="383D0896" w14:textId="780903CD" w:rsidR="00487D15" w:rsidRDefault="009D54E1" w:rsidP="00487D15"><w:pPr><w:pStyle w:val="a3"/><w:numPr><w:ilvl w:val="1"/><w:numId w:val="2"/></w:numPr></w:pPr><w:proofErr w:type="spellStart"/><w:r><w:t>Yafei</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> Liu, </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>Jiale</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> Wang</w:t></w:r><w:r w:rsidR="00487D15"><w:t xml:space="preserve"> – </w:t></w:r><w:r><w:t>Build Website</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$null = $taskRange.InsertXML($splitXml)
